# TC03_Canine_Filter_Diagnosis-Lymphoma.xlsx
# "Test cases with the Bioblank and Diagnosis filters"
#
# - Replace the "StatQuery" column (C2:C5) with the new Bioblank-style
#   counts query (Programs / Studies / Cases / Samples / Case Files / Study Files).
# - Replace the StudyFilesTab "query" cell (B5) with the new Diagnosis-filter
#   query for the Study Files tab (reworked WHERE clause + order/limit footer).
# - Row heights are nudged to match the new wrapped-text extents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bioblankQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
  WHERE diag.disease_term IN ['Lymphoma'] 
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$studyFilesQuery = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
WHERE diag.disease_term IN ['Lymphoma']
MATCH (c)<--(demo:demographic)
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

# Column C ("StatQuery") on every tab row now runs the Bioblank counts query.
$ws.Range("C2").Value = $bioblankQuery
$ws.Range("C3").Value = $bioblankQuery
$ws.Range("C4").Value = $bioblankQuery
$ws.Range("C5").Value = $bioblankQuery

# StudyFilesTab's "query" cell gets the reworked Diagnosis-filter query.
$ws.Range("B5").Value = $studyFilesQuery

# Match the new wrapped-text row heights.
$ws.Rows.Item(2).RowHeight = 333.5
$ws.Rows.Item(3).RowHeight = 232
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 377
